# Append " *" to a set of specific numeric result cells in the raw-data
# table (tolerance / barcode-search flag). Each search string is the
# exact original cell text; MatchWholeWord is used so that short values
# like "3.67" don't accidentally match inside longer ones like "3.671"
# or "3.678". wdReplaceAll (2) is used per call since every search
# string is either unique in the document or (for "3.671") all of its
# occurrences need the identical change.

$d = $word.ActiveDocument

$targets = @(
    "3.324",
    "3.372",
    "3.636",
    "3.197",
    "3.359",
    "3.375",
    "3.316",
    "3.33",
    "3.671",
    "3.164",
    "3.667",
    "3.669",
    "3.67"
)

foreach ($old in $targets) {
    $new = "$old *"
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}
